# "Can attach pdf files"
# Give the first two payslip rows their own per-person PDF location instead
# of the shared "dirnume1/payslip.pdf" location used by everybody.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "pers1/payslip.pdf"
$ws.Range("C3").Value = "pers2/payslip.pdf"

# View state left behind after making/reviewing the edit.
$excel.ActiveWindow.Zoom = 140
$ws.Range("B13").Select()
